$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Re-write header row (same text, but this moves them to "new" shared-string slots
# when the underlying strings table gets rebuilt on save).
$ws.Range("A1").Value = "run_num"
$ws.Range("B1").Value = "block_num"
$ws.Range("C1").Value = "start_time"
$ws.Range("D1").Value = "play_duration"
$ws.Range("E1").Value = "ear"
$ws.Range("F1").Value = "hand"

# Reset the four placeholder data rows (2-5) back to zeroed defaults and clear
# the ear/hand columns for those rows.
$ws.Range("A2:D5").Value = 0
$ws.Range("E2:F5").ClearContents()

# Column E (ear) got narrower.
$ws.Columns.Item(5).ColumnWidth = 3.1666666666666665
